$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate()

# Insert a new column before column N (14th column), shifting "Late" and
# "Outstanding" (and the blank spacer column) one column to the right.
$ws.Columns("N").Insert()

# Give the newly inserted column (N) the same width the author set for it.
# (ColumnWidth is offset from the stored OOXML width by ~0.8333 in this
# runtime, so request a ColumnWidth that serializes to width="10".)
$ws.Columns("N").ColumnWidth = 9.166666666666666

# Update the selection on the sheet to match the edited state.
$ws.Range("S6").Select()
